$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend rows 15-17 in column A with the same style as existing data rows (copy format from A14)
$ws.Range("A14").Copy($ws.Range("A15")) | Out-Null
$ws.Range("A14").Copy($ws.Range("A16")) | Out-Null
$ws.Range("A14").Copy($ws.Range("A17")) | Out-Null

# Row 2
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "youngstown-state"
$ws.Cells.Item(2,3).Value = "2017-18"
$ws.Cells.Item(2,4).Value = "Cameron Morse"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = "G"
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = "74.0"
$ws.Cells.Item(2,7).Style = "Normal"
$ws.Cells.Item(2,8).NumberFormat = "@"
$ws.Cells.Item(2,8).Value = "180"
$ws.Cells.Item(2,8).Style = "Normal"

# Row 3
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "youngstown-state"
$ws.Cells.Item(3,3).Value = "2017-18"
$ws.Cells.Item(3,4).Value = "Braun Hartfield"
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = "G"
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = "76.0"
$ws.Cells.Item(3,7).Style = "Normal"
$ws.Cells.Item(3,8).NumberFormat = "@"
$ws.Cells.Item(3,8).Value = "185"
$ws.Cells.Item(3,8).Style = "Normal"

# Row 4
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "youngstown-state"
$ws.Cells.Item(4,3).Value = "2017-18"
$ws.Cells.Item(4,4).Value = "Garrett Covington"
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = "G"
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = "77.0"
$ws.Cells.Item(4,7).Style = "Normal"
$ws.Cells.Item(4,8).NumberFormat = "@"
$ws.Cells.Item(4,8).Value = "195"
$ws.Cells.Item(4,8).Style = "Normal"

# Row 5
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "youngstown-state"
$ws.Cells.Item(5,3).Value = "2017-18"
$ws.Cells.Item(5,4).Value = "Naz Bohannon"
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = "F"
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = "78.0"
$ws.Cells.Item(5,7).Style = "Normal"
$ws.Cells.Item(5,8).NumberFormat = "@"
$ws.Cells.Item(5,8).Value = "230"
$ws.Cells.Item(5,8).Style = "Normal"

# Row 6
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "youngstown-state"
$ws.Cells.Item(6,3).Value = "2017-18"
$ws.Cells.Item(6,4).Value = "Tyree Robinson"
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = "F"
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = "77.0"
$ws.Cells.Item(6,7).Style = "Normal"
$ws.Cells.Item(6,8).NumberFormat = "@"
$ws.Cells.Item(6,8).Value = "220"
$ws.Cells.Item(6,8).Style = "Normal"

# Row 7
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "youngstown-state"
$ws.Cells.Item(7,3).Value = "2017-18"
$ws.Cells.Item(7,4).Value = "Jaylen Benton"
$ws.Cells.Item(7,5).Value = 0
$ws.Cells.Item(7,6).Value = "G"
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = "75.0"
$ws.Cells.Item(7,7).Style = "Normal"
$ws.Cells.Item(7,8).NumberFormat = "@"
$ws.Cells.Item(7,8).Value = "180"
$ws.Cells.Item(7,8).Style = "Normal"

# Row 8
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "youngstown-state"
$ws.Cells.Item(8,3).Value = "2017-18"
$ws.Cells.Item(8,4).Value = "Devin Haygood"
$ws.Cells.Item(8,5).Value = 0
$ws.Cells.Item(8,6).Value = "F"
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = "79.0"
$ws.Cells.Item(8,7).Style = "Normal"
$ws.Cells.Item(8,8).NumberFormat = "@"
$ws.Cells.Item(8,8).Value = "175"
$ws.Cells.Item(8,8).Style = "Normal"

# Row 9
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "youngstown-state"
$ws.Cells.Item(9,3).Value = "2017-18"
$ws.Cells.Item(9,4).Value = "Jeremiah Ferguson"
$ws.Cells.Item(9,5).Value = 0
$ws.Cells.Item(9,6).Value = "G"
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = "74.0"
$ws.Cells.Item(9,7).Style = "Normal"
$ws.Cells.Item(9,8).NumberFormat = "@"
$ws.Cells.Item(9,8).Value = "180"
$ws.Cells.Item(9,8).Style = "Normal"

# Row 10
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "youngstown-state"
$ws.Cells.Item(10,3).Value = "2017-18"
$ws.Cells.Item(10,4).Value = "Michael Akuchie"
$ws.Cells.Item(10,5).Value = 0
$ws.Cells.Item(10,6).Value = "F"
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = "80.0"
$ws.Cells.Item(10,7).Style = "Normal"
$ws.Cells.Item(10,8).NumberFormat = "@"
$ws.Cells.Item(10,8).Value = "215"
$ws.Cells.Item(10,8).Style = "Normal"

# Row 11
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "youngstown-state"
$ws.Cells.Item(11,3).Value = "2017-18"
$ws.Cells.Item(11,4).Value = "Noe Anabir"
$ws.Cells.Item(11,5).Value = 0
$ws.Cells.Item(11,6).Value = "F"
$ws.Cells.Item(11,7).NumberFormat = "@"
$ws.Cells.Item(11,7).Value = "80.0"
$ws.Cells.Item(11,7).Style = "Normal"
$ws.Cells.Item(11,8).NumberFormat = "@"
$ws.Cells.Item(11,8).Value = "226"
$ws.Cells.Item(11,8).Style = "Normal"

# Row 12
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "youngstown-state"
$ws.Cells.Item(12,3).Value = "2017-18"
$ws.Cells.Item(12,4).Value = "Francisco Santiago"
$ws.Cells.Item(12,5).Value = 0
$ws.Cells.Item(12,6).Value = "G"
$ws.Cells.Item(12,7).NumberFormat = "@"
$ws.Cells.Item(12,7).Value = "73.0"
$ws.Cells.Item(12,7).Style = "Normal"
$ws.Cells.Item(12,8).NumberFormat = "@"
$ws.Cells.Item(12,8).Value = "160"
$ws.Cells.Item(12,8).Style = "Normal"

# Row 13
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "youngstown-state"
$ws.Cells.Item(13,3).Value = "2017-18"
$ws.Cells.Item(13,4).Value = "Ryan Strollo"
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(13,6).Value = "G"
$ws.Cells.Item(13,7).NumberFormat = "@"
$ws.Cells.Item(13,7).Value = "74.0"
$ws.Cells.Item(13,7).Style = "Normal"
$ws.Cells.Item(13,8).NumberFormat = "@"
$ws.Cells.Item(13,8).Value = "190"
$ws.Cells.Item(13,8).Style = "Normal"

# Row 14
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "youngstown-state"
$ws.Cells.Item(14,3).Value = "2017-18"
$ws.Cells.Item(14,4).Value = "Jacob Brown"
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = "F"
$ws.Cells.Item(14,7).NumberFormat = "@"
$ws.Cells.Item(14,7).Value = "81.0"
$ws.Cells.Item(14,7).Style = "Normal"
$ws.Cells.Item(14,8).NumberFormat = "@"
$ws.Cells.Item(14,8).Value = "195"
$ws.Cells.Item(14,8).Style = "Normal"

# Row 15
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "youngstown-state"
$ws.Cells.Item(15,3).Value = "2017-18"
$ws.Cells.Item(15,4).Value = "Dan Ritter"
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = "G"
$ws.Cells.Item(15,7).NumberFormat = "@"
$ws.Cells.Item(15,7).Value = "75.0"
$ws.Cells.Item(15,7).Style = "Normal"
$ws.Cells.Item(15,8).NumberFormat = "@"
$ws.Cells.Item(15,8).Value = "176"
$ws.Cells.Item(15,8).Style = "Normal"

# Row 16
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "youngstown-state"
$ws.Cells.Item(16,3).Value = "2017-18"
$ws.Cells.Item(16,4).Value = "John Kirincic"
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = "G"
$ws.Cells.Item(16,7).NumberFormat = "@"
$ws.Cells.Item(16,7).Value = "73.0"
$ws.Cells.Item(16,7).Style = "Normal"
$ws.Cells.Item(16,8).NumberFormat = "@"
$ws.Cells.Item(16,8).Value = "175"
$ws.Cells.Item(16,8).Style = "Normal"

# Row 17
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "youngstown-state"
$ws.Cells.Item(17,3).Value = "2017-18"
$ws.Cells.Item(17,4).Value = "Alex Wilbourn"
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = "C"
$ws.Cells.Item(17,7).NumberFormat = "@"
$ws.Cells.Item(17,7).Value = "84.0"
$ws.Cells.Item(17,7).Style = "Normal"
$ws.Cells.Item(17,8).NumberFormat = "@"
$ws.Cells.Item(17,8).Value = "204"
$ws.Cells.Item(17,8).Style = "Normal"
